# Btc-Erbb2 LR-pairs sheet refresh: new TPM-derived NATMI output.
# The two rows whose "Target cluster" was "Neutrophils" no longer exist in the
# new run, so they are removed outright (old row 12, then old row 6 -- deleted
# high-to-low so row indices stay valid), and every remaining row/column is
# rewritten with the refreshed figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(12).Delete()
$ws.Rows(6).Delete()

# Row 2: FAPs -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Btc"
$ws.Range("C2").Value = "Erbb2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.06428666666666667
$ws.Range("H2").Value = 0.19286
$ws.Range("I2").Value = 0.005899579586891775
$ws.Range("J2").Value = 0.008823342375055644
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.138389
$ws.Range("N2").Value = 8.276778
$ws.Range("O2").Value = 0.3843663898945971
$ws.Range("P2").Value = 0.3371083667718477
$ws.Range("Q2").Value = 0.26604323418
$ws.Range("R2").Value = 1.59625940508
$ws.Range("S2").Value = 0.00226760010770945
$ws.Range("T2").Value = 0.002974422537523844

# Row 3: FAPs -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Btc"
$ws.Range("C3").Value = "Erbb2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.06428666666666667
$ws.Range("H3").Value = 0.19286
$ws.Range("I3").Value = 0.005899579586891775
$ws.Range("J3").Value = 0.008823342375055644
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.013819333333333
$ws.Range("N3").Value = 9.041457999999999
$ws.Range("O3").Value = 0.2799183104700584
$ws.Range("P3").Value = 0.3682533395986042
$ws.Range("Q3").Value = 0.1937483988755555
$ws.Range("R3").Value = 1.74373558988
$ws.Range("S3").Value = 0.001651400350446391
$ws.Range("T3").Value = 0.003249225296036121

# Row 4: FAPs -> Inflammatory-Mac
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Btc"
$ws.Range("C4").Value = "Erbb2"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.06428666666666667
$ws.Range("H4").Value = 0.19286
$ws.Range("I4").Value = 0.005899579586891775
$ws.Range("J4").Value = 0.008823342375055644
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.004857333333333333
$ws.Range("N4").Value = 0.014572
$ws.Range("O4").Value = 0.0004511406921505018
$ws.Range("P4").Value = 0.0005935091071186595
$ws.Range("Q4").Value = 0.0003122617688888889
$ws.Range("R4").Value = 0.00281035592
$ws.Range("S4").Value = 0.000002661540418227327
$ws.Range("T4").Value = 0.000005236734054821508

# Row 5: FAPs -> MuSCs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Btc"
$ws.Range("C5").Value = "Erbb2"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.06428666666666667
$ws.Range("H5").Value = 0.19286
$ws.Range("I5").Value = 0.005899579586891775
$ws.Range("J5").Value = 0.008823342375055644
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.609679
$ws.Range("N5").Value = 7.219358
$ws.Range("O5").Value = 0.3352607224473918
$ws.Range("P5").Value = 0.2940402635568179
$ws.Range("Q5").Value = 0.2320542306466667
$ws.Range("R5").Value = 1.39232538388
$ws.Range("S5").Value = 0.001977897314437222
$ws.Range("T5").Value = 0.002594417917413402

# Row 6: FAPs -> Resolving-Mac
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Btc"
$ws.Range("C6").Value = "Erbb2"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.06428666666666667
$ws.Range("H6").Value = 0.19286
$ws.Range("I6").Value = 0.005899579586891775
$ws.Range("J6").Value = 0.008823342375055644
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.000037
$ws.Range("N6").Value = 0.000111
$ws.Range("O6").Value = 0.000003436495802134621
$ws.Range("P6").Value = 0.000004520965611458359
$ws.Range("Q6").Value = 0.000002378606666666667
$ws.Range("R6").Value = 0.00002140746
$ws.Range("S6").Value = 0.00000002027388048471269
$ws.Range("T6").Value = 0.00000003989002745574989

# Row 7: MuSCs -> ECs
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Btc"
$ws.Range("C7").Value = "Erbb2"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 10.832535
$ws.Range("H7").Value = 21.66507
$ws.Range("I7").Value = 0.9941004204131083
$ws.Range("J7").Value = 0.9911766576249443
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.138389
$ws.Range("N7").Value = 8.276778
$ws.Range("O7").Value = 0.3843663898945971
$ws.Range("P7").Value = 0.3371083667718477
$ws.Range("Q7").Value = 44.829243686115
$ws.Range("R7").Value = 179.31697474446
$ws.Range("S7").Value = 0.3820987897868877
$ws.Range("T7").Value = 0.3341339442343239

# Row 8: MuSCs -> FAPs
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Btc"
$ws.Range("C8").Value = "Erbb2"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 10.832535
$ws.Range("H8").Value = 21.66507
$ws.Range("I8").Value = 0.9941004204131083
$ws.Range("J8").Value = 0.9911766576249443
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.013819333333333
$ws.Range("N8").Value = 9.041457999999999
$ws.Range("O8").Value = 0.2799183104700584
$ws.Range("P8").Value = 0.3682533395986042
$ws.Range("Q8").Value = 32.64730341201
$ws.Range("R8").Value = 195.88382047206
$ws.Range("S8").Value = 0.278266910119612
$ws.Range("T8").Value = 0.3650041143025681

# Row 9: MuSCs -> Inflammatory-Mac
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Btc"
$ws.Range("C9").Value = "Erbb2"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 10.832535
$ws.Range("H9").Value = 21.66507
$ws.Range("I9").Value = 0.9941004204131083
$ws.Range("J9").Value = 0.9911766576249443
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.004857333333333333
$ws.Range("N9").Value = 0.014572
$ws.Range("O9").Value = 0.0004511406921505018
$ws.Range("P9").Value = 0.0005935091071186595
$ws.Range("Q9").Value = 0.05261723334
$ws.Range("R9").Value = 0.31570340004
$ws.Range("S9").Value = 0.0004484791517322745
$ws.Range("T9").Value = 0.000588272373063838

# Row 10: MuSCs -> MuSCs
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Btc"
$ws.Range("C10").Value = "Erbb2"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 10.832535
$ws.Range("H10").Value = 21.66507
$ws.Range("I10").Value = 0.9941004204131083
$ws.Range("J10").Value = 0.9911766576249443
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.609679
$ws.Range("N10").Value = 7.219358
$ws.Range("O10").Value = 0.3352607224473918
$ws.Range("P10").Value = 0.2940402635568179
$ws.Range("Q10").Value = 39.101974106265
$ws.Range("R10").Value = 156.40789642506
$ws.Range("S10").Value = 0.3332828251329546
$ws.Range("T10").Value = 0.2914458456394046

# Row 11: MuSCs -> Resolving-Mac
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Btc"
$ws.Range("C11").Value = "Erbb2"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 10.832535
$ws.Range("H11").Value = 21.66507
$ws.Range("I11").Value = 0.9941004204131083
$ws.Range("J11").Value = 0.9911766576249443
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.000037
$ws.Range("N11").Value = 0.000111
$ws.Range("O11").Value = 0.000003436495802134621
$ws.Range("P11").Value = 0.000004520965611458359
$ws.Range("Q11").Value = 0.000400803795
$ws.Range("R11").Value = 0.00240482277
$ws.Range("S11").Value = 0.000003416221921649909
$ws.Range("T11").Value = 0.000004481075584002609
